$wb = $excel.ActiveWorkbook

# The "Belgium" sheet is the closest template for the new "Czech" sheet:
# same layout/merges/styles, just different market name + ticket number.
$belgium = $wb.Worksheets.Item("Belgium")

# Duplicate Belgium (keeps formatting, merged cells, column widths, styles)
# and place the copy at the end of the tab strip.
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$czech = $wb.Worksheets.Item($wb.Worksheets.Count)
$czech.Name = "Czech"

# Update the market name and ticket reference on the new sheet.
$czech.Range("B2").Value = "Czech Market"
$czech.Range("B4").Value = "NGC-3477/T1732"

# Restore Belgium's selection to the full used range (no longer the active tab)...
[void]$belgium.Range("A1:D15").Select()

# ...and make Czech the active tab with B4 selected.
$czech.Activate()
[void]$czech.Range("B4").Select()
